$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- HBAS Phase II Interim 1A by Will & 2A ---
# Add a new error-message entry "beas_wo_1002" (finish-goods serial number
# with no matching raw-material serial number) to both the Chinese (CH)
# and English (E) language blocks, right after the existing "beas_wo_1001"
# row in each block.

# Chinese (CH) block: insert new row 3, directly after row 2 (beas_wo_1001)
$ws.Rows("3:3").Insert()
$ws.Range("A3").Value = "CH  "
$ws.Range("B3").Value = "NULL"
$ws.Range("C3").Value = "beas_wo_1002"
$ws.Range("E3").Value = "成品收货的序列号<dw_1.item.serialnumber.value>，没有找到对应的原料序列号"

# English (E) block: insert new row 12, directly after row 11 (beas_wo_1001)
$ws.Rows("12:12").Insert()
$ws.Range("A12").Value = "E   "
$ws.Range("B12").Value = "NULL"
$ws.Range("C12").Value = "beas_wo_1002"
$ws.Range("E12").Value = "The serial number of finish goods <dw_1.item.serialnumber.value>, did not find the match raw material's serial number"

# Match the author's final selection state
$ws.Range("E21").Select()
